$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Junho (sheet8) - atualizacao da planilha de gastos
# ---------------------------------------------------------------------------
$wsJunho = $wb.Worksheets.Item("Junho")

# Linha 5: gasto "Gabinete PC" passa da coluna C (nao gasto) para a coluna B
# (gasto), com valor atualizado
$wsJunho.Range("C5").Clear()
$wsJunho.Range("B5").Value = 339.49

# Linha 6: gasto "Memoria 4GB 1333Mhz" passa da coluna C para a coluna B
$wsJunho.Range("C6").Clear()
$wsJunho.Range("B6").Value = 123

# Linha 7: gasto "Mouse Bluetooth" passa da coluna C para a coluna B
$wsJunho.Range("C7").Clear()
$wsJunho.Range("B7").Value = 81.72

# Linha 10: novo gasto "Cooler Led"
$wsJunho.Range("A10").Value = "Cooler Led"
$wsJunho.Range("B10").Value = 48.46

# Linha 11: nova linha "PosEAD" (mesmo padrao usado nas outras abas)
$wsJunho.Range("A11").Value = "PosEAD"
$wsJunho.Range("A11").Font.Name = "Calibri"
$wsJunho.Range("A11").Font.Size = 11
$wsJunho.Range("A11").Font.Color = 0
$wsJunho.Rows.Item(11).RowHeight = 15
$wsJunho.Range("B11").Value = 210.83

# ---------------------------------------------------------------------------
# Maio (sheet7) - apenas a selecao ativa foi alterada
# ---------------------------------------------------------------------------
$wsMaio = $wb.Worksheets.Item("Maio")
$wsMaio.Activate()
$wsMaio.Range("A5").Select()

# ---------------------------------------------------------------------------
# Volta para Junho (aba ativa original) e atualiza a celula selecionada
# ---------------------------------------------------------------------------
$wsJunho.Activate()
$wsJunho.Range("B12").Select()
